$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = "current"
$ws.Range("D20").Value = 0.01
$ws.Range("B21").Formula = "=(B20-B22)/D20"
$ws.Range("B31").Value = 100

$ws.Range("D20").Select()
